# Sync attendance_reports: swap "Recorded By" (column G) name order
# from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# wherever the old value is found.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
